$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update report title (page_1_report_title) and firm name (page_1_firm_name)
$ws.Range("C4").Value = "Quarterly Portfolio Report"
$ws.Range("C5").Value = "Gaard Capital LLC"

# Match the resulting selection/active cell from the edit session
$ws.Range("E11").Select()
